# Edit script: regenerate "Programmes pour les enfants" child-program input cells
# (commit: "Regeneration of en/fr/es country data books after child program paras updated")
#
# The source data table on "Programmes pour les enfants" (columns D:H, rows 2-53)
# holds literal input percentages; columns F:H (and a few D cells) were
# recomputed upstream and need to be written back here. Every dependent
# formula cell below row 53 (rows 57-163) references these inputs via
# shared formulas (e.g. "=D18*0.9", "=D2*0.9", etc.) and will recalculate
# automatically once the inputs change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Programmes pour les enfants")

$edits = @(
    @{Cell="F2"; Value=0.39473684210526322}
    @{Cell="G2"; Value=0.39473684210526322}
    @{Cell="H2"; Value=0.39473684210526322}
    @{Cell="F3"; Value=0.30769230769230765}
    @{Cell="G3"; Value=0.30769230769230765}
    @{Cell="H3"; Value=0.30769230769230765}
    @{Cell="F18"; Value=0.7}
    @{Cell="F20"; Value=0.84}
    @{Cell="D21"; Value=0.28260869565217389}
    @{Cell="F21"; Value=0}
    @{Cell="F22"; Value=0}
    @{Cell="D23"; Value=0.28260869565217389}
    @{Cell="F23"; Value=0}
    @{Cell="F24"; Value=0}
    @{Cell="D25"; Value=0.28260869565217389}
    @{Cell="F25"; Value=0}
    @{Cell="F26"; Value=0}
    @{Cell="F27"; Value=1}
    @{Cell="F28"; Value=0}
    @{Cell="F29"; Value=0}
    @{Cell="F30"; Value=1}
    @{Cell="F31"; Value=0}
    @{Cell="F32"; Value=0}
    @{Cell="F33"; Value=1}
    @{Cell="F34"; Value=0}
    @{Cell="F35"; Value=0}
    @{Cell="F36"; Value=1}
    @{Cell="F37"; Value=0}
    @{Cell="F38"; Value=0}
    @{Cell="F39"; Value=1}
    @{Cell="F40"; Value=0}
    @{Cell="F41"; Value=0}
    @{Cell="F42"; Value=0.3}
    @{Cell="F43"; Value=0.5}
    @{Cell="F44"; Value=0.65}
    @{Cell="F45"; Value=0.3}
    @{Cell="F46"; Value=0.49}
    @{Cell="F47"; Value=0.52}
    @{Cell="F48"; Value=0.88}
    @{Cell="D49"; Value=0.78409090909090906}
    @{Cell="E49"; Value=0.78409090909090906}
    @{Cell="F49"; Value=0.78409090909090906}
    @{Cell="G49"; Value=0.78409090909090906}
    @{Cell="H49"; Value=0.78409090909090906}
    @{Cell="D50"; Value=0.88372093023255816}
    @{Cell="E50"; Value=0.88372093023255816}
    @{Cell="F50"; Value=0.88372093023255816}
    @{Cell="G50"; Value=0.88372093023255816}
    @{Cell="H50"; Value=0.88372093023255816}
    @{Cell="F51"; Value=0.86}
    @{Cell="F52"; Value=0}
    @{Cell="F53"; Value=0}
)

foreach ($edit in $edits) {
    $ws.Range($edit.Cell).Value = $edit.Value
}

# View-state bookkeeping to mirror the author's session: the
# "Programmes pour les enfants" selection moved to D2:H53, and the workbook's
# active sheet moved from "Donnees pop de l'annee de ref" to
# "Dependances du programme" (activated last so it ends up the active tab).
$ws.Activate() | Out-Null
$ws.Range("D2:H53").Select() | Out-Null

$wsDep = $wb.Worksheets.Item("Dépendances du programme")
$wsDep.Activate() | Out-Null
$wsDep.Range("A19").Select() | Out-Null

$wb.Application.Calculate()
